$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" values; regenerate per commit ("use K instead of Strike#").
$kValues = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 0
    6  = 2
    7  = 1
    8  = 2
    9  = 0
    10 = 2
    11 = 1
    12 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
